# Restored from revision of admin on 05/18/2020 07:59:41 AM.TEST Author: admin. Type: SAVE.
# Main.xlsx / "Rules" sheet: cell C10 ("Integer min" for rule R30) changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Rules")

$ws.Range("C10").Value = 1
